$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $row, $col, $value) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-CellText $ws 2 4 '28.224.61'
Set-CellText $ws 2 5 '  +0.18%  '
Set-CellText $ws 3 4 '1.855.35'
Set-CellText $ws 3 5 '  -0.69%  '
Set-CellText $ws 4 5 '  +0.21%  '
Set-CellText $ws 5 4 '329.50'
Set-CellText $ws 5 5 '  -2.19%  '
Set-CellText $ws 6 4 '1.001'
Set-CellText $ws 6 5 '  +0.11%  '
Set-CellText $ws 7 4 '0.4545'
Set-CellText $ws 8 4 '0.3921'
Set-CellText $ws 8 5 '  -0.07%  '
Set-CellText $ws 9 4 '47.51'
Set-CellText $ws 9 5 '  +1.51%  '
Set-CellText $ws 10 4 '0.07793'
Set-CellText $ws 10 5 '  -2.29%  '
Set-CellText $ws 11 4 '0.9804'
Set-CellText $ws 11 5 '  -2.19%  '
Set-CellText $ws 12 4 '21.24'
Set-CellText $ws 12 5 '  -2.28%  '
Set-CellText $ws 13 4 '1.879.61'
Set-CellText $ws 13 5 '  +0.59%  '
Set-CellText $ws 14 4 '5.789'
Set-CellText $ws 14 5 '  -3.22%  '
Set-CellText $ws 15 4 '6.942'
Set-CellText $ws 15 5 '  -4.49%  '
Set-CellText $ws 16 4 '1.003'
Set-CellText $ws 16 5 '  +0.14%  '
Set-CellText $ws 17 4 '87.80'
Set-CellText $ws 17 5 '  -3.94%  '
Set-CellText $ws 18 4 '0.06523'
Set-CellText $ws 18 5 '  -1.10%  '
Set-CellText $ws 19 4 '0.00001015'
Set-CellText $ws 19 5 '  -2.75%  '
Set-CellText $ws 20 4 '17.04'
Set-CellText $ws 20 5 '  -4.06%  '
Set-CellText $ws 21 5 '  +0.16%  '
Set-CellText $ws 22 4 '28.228.69'
Set-CellText $ws 22 5 '  +0.20%  '
Set-CellText $ws 23 4 '5.286'
Set-CellText $ws 23 5 '  -2.83%  '
Set-CellText $ws 24 4 '10.68'
Set-CellText $ws 24 5 '  -3.37%  '
Set-CellText $ws 25 4 '2.255'
Set-CellText $ws 25 5 '  -1.70%  '
Set-CellText $ws 26 4 '2.082.61'
Set-CellText $ws 26 5 '  +0.13%  '
Set-CellText $ws 27 4 '156.92'
Set-CellText $ws 27 5 '  -1.38%  '
Set-CellText $ws 28 4 '19.15'
Set-CellText $ws 28 5 '  -3.27%  '
Set-CellText $ws 29 4 '2.040'
Set-CellText $ws 29 5 '  -4.09%  '
Set-CellText $ws 30 4 '5.254'
Set-CellText $ws 30 5 '  -4.27%  '
Set-CellText $ws 31 4 '116.09'
Set-CellText $ws 31 5 '  -3.02%  '
Set-CellText $ws 32 4 '0.9377'
Set-CellText $ws 32 5 '  -3.70%  '
Set-CellText $ws 33 4 '0.09260'
Set-CellText $ws 33 5 '  -2.43%  '
Set-CellText $ws 34 4 '3.605'
Set-CellText $ws 34 5 '  +0.89%  '
Set-CellText $ws 35 4 '1.373'
Set-CellText $ws 35 5 '  -0.48%  '
Set-CellText $ws 36 4 '5.182'
Set-CellText $ws 36 5 '  -3.01%  '
Set-CellText $ws 37 4 '0.06004'
Set-CellText $ws 37 5 '  -1.50%  '
Set-CellText $ws 38 4 '0.02191'
Set-CellText $ws 38 5 '  -3.53%  '
Set-CellText $ws 39 4 '8.167'
Set-CellText $ws 39 5 '  -2.80%  '
Set-CellText $ws 40 4 '1.160'
Set-CellText $ws 40 5 '  -1.17%  '
Set-CellText $ws 41 5 '  +0.22%  '
Set-CellText $ws 42 4 '0.5651'
Set-CellText $ws 42 5 '  -5.55%  '
Set-CellText $ws 43 4 '9.950'
Set-CellText $ws 43 5 '  -3.69%  '
Set-CellText $ws 44 4 '0.1788'
Set-CellText $ws 44 5 '  -5.05%  '
Set-CellText $ws 45 4 '1.245'
Set-CellText $ws 45 5 '  -2.40%  '
Set-CellText $ws 46 4 '2.286'
Set-CellText $ws 46 5 '  +21.03%  '
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText $ws 47 4 '11.78'
Set-CellText $ws 47 5 '  -2.91%  '
$ws.Cells.Item(48, 2).Value = 'Decentraland'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-CellText $ws 48 4 '0.5354'
Set-CellText $ws 48 5 '  -4.86%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws 49 4 '0.07168'
Set-CellText $ws 49 5 '  +4.50%  '
Set-CellText $ws 50 4 '1.858'
Set-CellText $ws 50 5 '  -6.13%  '
Set-CellText $ws 51 4 '109.44'
Set-CellText $ws 51 5 '  -1.91%  '
